# Fruta / hortaliza, semanal
# Insert a new weekly record at row 42 (pushing the existing rows 42..89 down
# to 43..90) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 42; this shifts rows 42-89
# down to 43-90 and extends the used range accordingly.
$ws.Rows.Item(42).Insert()

# Fill in the new row 42 with the new observation.
$ws.Cells.Item(42, 1).Value = 3
$ws.Cells.Item(42, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(42, 3).Value = "Coquimbo"
$ws.Cells.Item(42, 4).Value = 44586
$ws.Cells.Item(42, 5).Value = 5
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100107
$ws.Cells.Item(42, 8).Value = "Otros"
$ws.Cells.Item(42, 9).Value = 100107011
$ws.Cells.Item(42, 10).Value = "Tuna"
$ws.Cells.Item(42, 11).Value = "Sin especificar"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 120
$ws.Cells.Item(42, 14).Value = 19000
$ws.Cells.Item(42, 15).Value = 20000
$ws.Cells.Item(42, 16).Value = 19500
$ws.Cells.Item(42, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(42, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(42, 19).Value = 975
$ws.Cells.Item(42, 20).Value = 20
